$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 (sheet1) ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 14:58:38"
$ws1.Range("A3").Value = "Total filas: 200"
$ws1.Range("C15").Value = "225_GOMEZ"
$ws1.Range("C16").Value = "215A_EL PATO"
$ws1.Range("A35").Value = "07:26:49"
$ws1.Range("C35").Value = "16_SANTA ANA"
$ws1.Range("D35").Value = 34
$ws1.Range("A36").Value = "07:51:40"
$ws1.Range("C36").Value = "17_ROMERO"
$ws1.Range("D36").Value = 9
$ws1.Range("A43").Value = "08:14:55"
$ws1.Range("C43").Value = "11_ETCHEVERRY"
$ws1.Range("D43").Value = 14
$ws1.Range("A44").Value = "06:58:58"
$ws1.Range("C44").Value = "15_ABASTO"
$ws1.Range("D44").Value = 90
$ws1.Range("A71").Value = "08:49:06"
$ws1.Range("C71").Value = "23_HERNANDEZ"
$ws1.Range("D71").Value = 42
$ws1.Range("A72").Value = "08:14:55"
$ws1.Range("C72").Value = "16_SANTA ANA"
$ws1.Range("D72").Value = 77
$ws1.Range("C90").Value = "14_ABASTO"
$ws1.Range("C91").Value = "15_ABASTO"
$ws1.Range("A122").Value = "11:38:09"
$ws1.Range("C122").Value = "15_ABASTO"
$ws1.Range("D122").Value = 32
$ws1.Range("A123").Value = "11:56:32"
$ws1.Range("C123").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D123").Value = 14
$ws1.Range("C133").Value = "23_HERNANDEZ"
$ws1.Range("C134").Value = "27_EL RETIRO"
$ws1.Range("A136").Value = "12:18:38"
$ws1.Range("C136").Value = "16_SANTA ANA"
$ws1.Range("D136").Value = 30
$ws1.Range("C137").Value = "15X38_ABASTO"
$ws1.Range("A138").Value = "12:43:13"
$ws1.Range("C138").Value = "14_ABASTO"
$ws1.Range("D138").Value = 5
$ws1.Range("A140").Value = "12:43:13"
$ws1.Range("C140").Value = "11_ETCHEVERRY"
$ws1.Range("D140").Value = 20
$ws1.Range("A141").Value = "12:18:38"
$ws1.Range("C141").Value = "215C_EL PATO"
$ws1.Range("D141").Value = 45
$ws1.Range("A148").Value = "12:43:13"
$ws1.Range("C148").Value = "23_HERNANDEZ"
$ws1.Range("D148").Value = 42
$ws1.Range("A149").Value = "12:58:23"
$ws1.Range("C149").Value = "16_P MOR-SANTA ANA"
$ws1.Range("D149").Value = 27
$ws1.Range("A151").Value = "13:28:27"
$ws1.Range("C151").Value = "215A_EL PATO"
$ws1.Range("D151").Value = 5
$ws1.Range("A152").Value = "12:43:13"
$ws1.Range("C152").Value = "14_ABASTO"
$ws1.Range("D152").Value = 50
$ws1.Range("C162").Value = "27_EL RETIRO"
$ws1.Range("C163").Value = "11_ETCHEVERRY"
$ws1.Range("A177").Value = "14:58:38"
$ws1.Range("D177").Value = 4
$ws1.Range("A179").Value = "14:58:38"
$ws1.Range("D179").Value = 15
$ws1.Range("A180").Value = "14:58:38"
$ws1.Range("D180").Value = 18
$ws1.Range("A181").Value = "14:58:38"
$ws1.Range("D181").Value = 19
$ws1.Range("A184").Value = "14:58:38"
$ws1.Range("D184").Value = 35
$ws1.Range("A186").Value = "14:58:38"
$ws1.Range("D186").Value = 38
$ws1.Range("A187").Value = "14:58:38"
$ws1.Range("D187").Value = 43
$ws1.Range("A189").Value = "14:58:38"
$ws1.Range("D189").Value = 55
$ws1.Range("A190").Value = "14:58:38"
$ws1.Range("D190").Value = 55
$ws1.Range("A192").Value = "14:58:38"
$ws1.Range("B192").Value = "15:56"
$ws1.Range("D192").Value = 58
$ws1.Range("A193").Value = "14:17:13"
$ws1.Range("B193").Value = "15:57"
$ws1.Range("C193").Value = "27_EL RETIRO"
$ws1.Range("D193").Value = 100
$ws1.Range("A194").Value = "14:58:38"
$ws1.Range("B194").Value = "16:01"
$ws1.Range("C194").Value = "16_SANTA ANA"
$ws1.Range("D194").Value = 63
$ws1.Range("A195").Value = "14:58:38"
$ws1.Range("B195").Value = "16:05"
$ws1.Range("C195").Value = "14_ABASTO"
$ws1.Range("D195").Value = 67
$ws1.Range("A196").Value = "14:17:13"
$ws1.Range("B196").Value = "16:06"
$ws1.Range("C196").Value = "14_ABASTO"
$ws1.Range("D196").Value = 109
$ws1.Range("A197").Value = "14:58:38"
$ws1.Range("B197").Value = "16:14"
$ws1.Range("C197").Value = "17_ROMERO"
$ws1.Range("D197").Value = 76
$ws1.Range("B198").Value = "16:16"
$ws1.Range("C198").Value = "10_OLMOS"
$ws1.Range("D198").Value = 92
$ws1.Range("A199").Value = "14:58:38"
$ws1.Range("B199").Value = "16:17"
$ws1.Range("C199").Value = "10_OLMOS"
$ws1.Range("D199").Value = 79
$ws1.Range("A200").Value = "14:58:38"
$ws1.Range("B200").Value = "16:21"
$ws1.Range("C200").Value = "23_HERNANDEZ"
$ws1.Range("D200").Value = 83
$ws1.Range("E200").Value = "LP1912"
$ws1.Range("A201").Value = "14:58:38"
$ws1.Range("B201").Value = "16:34"
$ws1.Range("C201").Value = "83_ALUAR"
$ws1.Range("D201").Value = 96
$ws1.Range("E201").Value = "LP1912"
$ws1.Range("A202").Value = "14:44:54"
$ws1.Range("B202").Value = "16:40"
$ws1.Range("C202").Value = "225_GOMEZ"
$ws1.Range("D202").Value = 116
$ws1.Range("E202").Value = "LP1912"
$ws1.Range("A203").Value = "14:58:38"
$ws1.Range("B203").Value = "16:41"
$ws1.Range("C203").Value = "225_GOMEZ"
$ws1.Range("D203").Value = 103
$ws1.Range("E203").Value = "LP1912"
$ws1.Range("A204").Value = "14:58:38"
$ws1.Range("B204").Value = "16:46"
$ws1.Range("C204").Value = "17_ROMERO"
$ws1.Range("D204").Value = 108
$ws1.Range("E204").Value = "LP1912"
$ws1.Range("A205").Value = "14:58:38"
$ws1.Range("B205").Value = "16:53"
$ws1.Range("C205").Value = "11_ETCHEVERRY"
$ws1.Range("D205").Value = 115
$ws1.Range("E205").Value = "LP1912"

# ---- Sheet: LP1912-215 (sheet2) ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 14:58:38"
$ws2.Range("A37").Value = "14:58:38"
$ws2.Range("D37").Value = 35

# ---- Sheet: 6203-6173 (sheet3) ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 14:58:38"
$ws3.Range("A3").Value = "Total filas: 34"
$ws3.Range("C22").Value = "215B_LP-P MOR-1 Y 57"
$ws3.Range("C23").Value = "215A_LA PLATA"
$ws3.Range("A34").Value = "14:58:38"
$ws3.Range("D34").Value = 23
$ws3.Range("A36").Value = "14:58:38"
$ws3.Range("D36").Value = 63
$ws3.Range("A39").Value = "14:58:38"
$ws3.Range("B39").Value = "16:30"
$ws3.Range("C39").Value = "215B_LP-P MOR-40 Y 115"
$ws3.Range("D39").Value = 92
$ws3.Range("E39").Value = "L6173"
